$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial of 45178 for every data row
# (rows 2 through 205). Update it to 45179 to reflect the new "changed" date.
$ws.Range("C2:C205").Value = 45179
